# Add 2022-Q3 data
#  - insert a new top data row into the "总计" (summary) sheet
#  - insert a brand-new "2022-Q3" worksheet (positioned right after "总计",
#    before the existing "2022-Q1" sheet) holding the per-fund breakdown

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 总计 sheet: insert new row 2 (2022-Q3 summary), pushing the rest down
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Rows.Item(2).ClearFormats()
# match formatting of the surrounding data rows (style only lives on col A)
$summary.Range("A3").Copy($summary.Range("A2"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 8
$summary.Range("D2").Value = 0.78

# the A column is a sequential row index (0,1,2,...) - renumber the rows
# that got pushed down by the insert
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q3" sheet right after "总计" by duplicating the
#    "2022-Q1" sheet (same headers/column formatting), then overwrite it
#    with the 2022-Q3 fund holdings.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q1Sheet.Copy($null, $totalSheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Grow from 1 data row (copied from 2022-Q1) to 8 data rows, copying the
# formatting of row 2 down to rows 3-9.
for ($r = 3; $r -le 9; $r++) {
    $q3Sheet.Range("A2:H2").Copy($q3Sheet.Range("A" + $r + ":H" + $r))
}

# Columns B (fund code) and D:G (stored as text in this workbook, e.g. to
# keep leading zeros on fund codes) must stay text, not get coerced to
# numbers.
$q3Sheet.Range("B2:B9").NumberFormat = "@"
$q3Sheet.Range("D2:G9").NumberFormat = "@"

$data = @(
    @(0, "340008", "兴全有机增长混合",                 "23.13", "76.53", "3.04", "0.7032", 8),
    @(1, "000827", "广发中证百度百发策略100指数E",      "2.74",  "92.46", "1.05", "0.0288", 3),
    @(2, "010931", "国联安鑫元1个月持有期混合A",        "1.88",  "33.30", "1.23", "0.0231", 9),
    @(3, "000826", "广发中证百度百发策略100指数A",      "0.89",  "92.46", "1.05", "0.0093", 3),
    @(4, "005536", "渤海汇金量化成长混合",               "0.38",  "84.80", "1.63", "0.0062", 6),
    @(5, "010584", "渤海汇金新动能主题混合",             "0.27",  "91.08", "1.98", "0.0053", 5),
    @(6, "007808", "北信瑞丰量化优选灵活配置混合",       "0.17",  "78.75", "1.22", "0.0021", 2),
    @(7, "010932", "国联安鑫元1个月持有期混合C",        "0.10",  "33.30", "1.23", "0.0012", 9)
)

$row = 2
foreach ($item in $data) {
    $q3Sheet.Range("A" + $row).Value = $item[0]
    $q3Sheet.Range("B" + $row).Value = $item[1]
    $q3Sheet.Range("C" + $row).Value = $item[2]
    $q3Sheet.Range("D" + $row).Value = $item[3]
    $q3Sheet.Range("E" + $row).Value = $item[4]
    $q3Sheet.Range("F" + $row).Value = $item[5]
    $q3Sheet.Range("G" + $row).Value = $item[6]
    $q3Sheet.Range("H" + $row).Value = $item[7]
    $row = $row + 1
}
